$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.562.00'
$ws.Range("E2").Value = '  +2.62%  '

$ws.Range("D3").Value = '1.670.06'
$ws.Range("E3").Value = '  +2.08%  '

$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '239.09'
$ws.Range("E5").Value = '  +1.65%  '

$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '0.4779'
$ws.Range("E7").Value = '  +1.49%  '

$ws.Range("D8").Value = '0.2627'
$ws.Range("E8").Value = '  +3.12%  '

$ws.Range("D9").Value = '0.06172'
$ws.Range("E9").Value = '  +2.72%  '

$ws.Range("D10").Value = '1.669.89'
$ws.Range("E10").Value = '  +2.06%  '

$ws.Range("D11").Value = '0.06996'
$ws.Range("E11").Value = '  -2.23%  '

$ws.Range("E12").Value = '  +1.15%  '

$ws.Range("D13").Value = '0.5898'
$ws.Range("E13").Value = '  -3.76%  '

$ws.Range("D14").Value = '4.379'
$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("E15").Value = '  +3.89%  '

$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").Value = '25.554.40'
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("D19").Value = '0.000006769'
$ws.Range("E19").Value = '  +3.14%  '

$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("D21").Value = '1.884.18'
$ws.Range("E21").Value = '  +2.05%  '

$ws.Range("D22").Value = '4.444'
$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").Value = '8.751'
$ws.Range("E23").Value = '  +2.44%  '

$ws.Range("D24").Value = '5.276'
$ws.Range("E24").Value = '  +0.64%  '

$ws.Range("D25").Value = '136.62'
$ws.Range("E25").Value = '  +3.20%  '

$ws.Range("D26").Value = '15.03'
$ws.Range("E26").Value = '  +1.80%  '

$ws.Range("E27").Value = '  +1.31%  '

$ws.Range("E28").Value = '  +4.49%  '

$ws.Range("E29").Value = '  +2.09%  '

$ws.Range("D30").Value = "'3.960"
$ws.Range("E30").Value = '  +6.60%  '

$ws.Range("D31").Value = '0.07828'
$ws.Range("E31").Value = '  +1.10%  '

$ws.Range("D32").Value = '3.651'
$ws.Range("E32").Value = '  +3.46%  '

$ws.Range("D33").Value = '0.9992'
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").Value = '0.04224'
$ws.Range("E34").Value = '  -3.28%  '

$ws.Range("D35").Value = '2.623'
$ws.Range("E35").Value = '  +1.00%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6098'
$ws.Range("E36").Value = '  +5.23%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '0.9554'
$ws.Range("E37").Value = '  +4.52%  '

$ws.Range("D38").Value = '2.592'
$ws.Range("E38").Value = '  +2.20%  '

$ws.Range("D39").Value = '0.8613'
$ws.Range("E39").Value = '  +5.16%  '

$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("D41").Value = "'1.870"
$ws.Range("E41").Value = '  +4.59%  '

$ws.Range("D42").Value = '0.01478'
$ws.Range("E42").Value = '  -4.45%  '

$ws.Range("D43").Value = '96.61'
$ws.Range("E43").Value = '  -0.61%  '

$ws.Range("D44").Value = "'0.3770"
$ws.Range("E44").Value = '  +2.05%  '

$ws.Range("D45").Value = "'4.880"
$ws.Range("E45").Value = '  +3.11%  '

$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("D47").Value = '6.218'
$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("D48").Value = '0.05259'
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("E49").Value = '  +1.64%  '

$ws.Range("D50").Value = '7.403'
$ws.Range("E50").Value = '  +3.44%  '

$ws.Range("E51").Value = '  +0.13%  '
